# Load Factor speedups missing hazptr
#
# Fills in the "Hazard Pointer" (column F) speedup values on the
# "Sheet2" worksheet (the load_factor_test.txt data), which were
# previously missing, and refreshes the dependent "DCAS" (column F)
# values on the "Sheet1" worksheet (the 5% del / 5% ins load-factor
# chart source rows 59:65) to match the recomputed numbers.

$wb = $excel.ActiveWorkbook

# --- Sheet2: add the missing Hazard Pointer (column F) values ------------
$ws2 = $wb.Worksheets.Item("Sheet2")

$hazptrValues = @{
    3  = 0.28692299999999998
    4  = 0.44977200000000001
    5  = 0.55448299999999995
    6  = 0.70880799999999999
    7  = 0.71130199999999999
    8  = 0.717916
    9  = 0.71992500000000004
    11 = 0.55763499999999999
    12 = 0.73817500000000003
    13 = 1.215751
    14 = 1.302767
    15 = 1.2324790000000001
    16 = 1.2296640000000001
    17 = 1.205122
    19 = 0.37883600000000001
    20 = 0.59200399999999997
    21 = 0.874089
    22 = 0.90456199999999998
    23 = 0.91234599999999999
    24 = 0.91863600000000001
    25 = 0.92056300000000002
    27 = 0.41204000000000002
    28 = 0.51643099999999997
    29 = 0.75488599999999995
    30 = 0.98655000000000004
    31 = 0.98013499999999998
    32 = 0.98258400000000001
    33 = 0.95097600000000004
    35 = 0.45594899999999999
    36 = 0.75457700000000005
    37 = 1.0737639999999999
    38 = 1.1371290000000001
    39 = 1.0699689999999999
    40 = 1.0421609999999999
    41 = 1.0362199999999999
}

foreach ($row in $hazptrValues.Keys) {
    $ws2.Cells.Item($row, 6).Value = $hazptrValues[$row]
}

$ws2.Range("I42").Select() | Out-Null

# --- Sheet1: refresh the recomputed DCAS (column F) values, rows 59:65 ---
$ws1 = $wb.Worksheets.Item("Sheet1")

$dcasValues = @{
    59 = 0.43965500000000002
    60 = 0.79086900000000004
    61 = 1.42069
    62 = 1.837367
    63 = 1.9172210000000001
    64 = 1.903195
    65 = 1.9467669999999999
}

foreach ($row in $dcasValues.Keys) {
    $ws1.Cells.Item($row, 6).Value = $dcasValues[$row]
}

$ws1.Application.ActiveWindow.ScrollRow = 36
$ws1.Range("F69").Select() | Out-Null
